{"js": "// Week 7 - Meeting Report updates:\n//   1. \"Other Members Present: ... Gerjan Haxhia\" -> \"... Gerjan Haxhija\"\n//   2. \"Time: 14:40, April 11\" runs consolidated into a single run.\n\nconst body = context.document.body;\n\n// --- 1. Fix \"Haxhia\" -> \"Haxhija\" -----------------------------------------\nconst nameHits = body.search(\"Gerjan Haxhia\", { matchCase: true });\nnameHits.load(\"items\");\nawait context.sync();\n\nif (nameHits.items.length > 0) {\n  nameHits.items[0].insertText(\"Gerjan Haxhija\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2. Collapse the \"Time: 14:40, April 11\" runs into one ----------------\nconst timeHits = body.search(\"Time: 14:40, April 11\", { matchCase: true });\ntimeHits.load(\"items\");\nawait context.sync();\n\nif (timeHits.items.length > 0) {\n  timeHits.items[0].insertText(\"Time: 14:40, April 11\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Week 7 - Meeting Report updates:\n#   1. \"Other Members Present: ... Gerjan Haxhia\" -> \"... Gerjan Haxhija\"\n#   2. \"Time: 14:40, April 11\" runs consolidated into a single run.\n\n$d = $word.ActiveDocument\n\n# --- 1. Fix \"Haxhia\" -> \"Haxhija\" ------------------------------------------\n$find1 = $d.Content.Find\n$find1.Text = \"Gerjan Haxhia\"\n$find1.Replacement.Text = \"Gerjan Haxhija\"\n$find1.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n\n# --- 2. Collapse the \"Time: 14:40, April 11\" runs into one -----------------\n$find2 = $d.Content.Find\n$find2.Text = \"Time: 14:40, April 11\"\n$find2.Replacement.Text = \"Time: 14:40, April 11\"\n$find2.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2)\n"}
